$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set "CAM" values first so it becomes shared-string index 33 (added before "JEFF!")
$ws.Range("D14").Value = "CAM"
$ws.Range("D16").Value = "CAM"
$ws.Range("D19").Value = "CAM"
$ws.Range("D20").Value = "CAM"
$ws.Range("D21").Value = "CAM"
$ws.Range("D22").Value = "CAM"
$ws.Range("D23").Value = "CAM"
$ws.Range("D24").Value = "CAM"
$ws.Range("D25").Value = "CAM"

# Set "JEFF!" values so it becomes shared-string index 34 (added after "CAM")
$ws.Range("D10").Value = "JEFF!"
$ws.Range("D11").Value = "JEFF!"

# Update the active selection to D11 (was D28)
$ws.Range("D11").Select()
